$d = $word.ActiveDocument

$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$newPara = $r2.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertBefore("XXXXXXXXXXXXXXXXXXXXXXXXXXXXXXXX")

$p3b = $d.Paragraphs.Item(3)
$r3 = $p3b.Range
$newPara2 = $r3.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertBefore("Edits from the blue-dev1 branch")  # no trailing period yet

$p4c = $d.Paragraphs.Item(4)
Write-Output "p4: $($p4c.Range.Start)-$($p4c.Range.End)"

$d.Bookmarks("_GoBack").Delete()
$endPos = $p4c.Range.End - 1   # 126, cursed position, but let's put bookmark a bit before that and grow it
$safePos = $endPos - 3
Write-Output "safePos=$safePos endPos=$endPos"
$d.Bookmarks.Add("_GoBack", $d.Range($safePos, $safePos))
$b = $d.Bookmarks.Item("_GoBack")
Write-Output "bm at safe pos: $($b.Start)-$($b.End)"

# now insert the missing chars ("nch") was already present; let's instead test growing bookmark by inserting text right at its position from outside
$ins = $d.Range($safePos, $safePos)
$ins.InsertAfter("Q")
$b2 = $d.Bookmarks.Item("_GoBack")
Write-Output "bm after insertAfter Q at same pos: $($b2.Start)-$($b2.End)"
